$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("G1").Value = "SLA (g/cm2)"
$ws.Range("J1").Value = "phi_total_no_yield"

# Update recalculated phi_total values for rows 4 and 11 (columns I and J)
$ws.Range("I4").Value = 1.29389041132425
$ws.Range("J4").Value = 1.29389041132425

$ws.Range("I11").Value = 1.35850618518056
$ws.Range("J11").Value = 1.35850618518056
